$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row (F1:H1) with new column names, copying the header
# formatting (bold font, thin border, centered/top alignment) used by A1:E1.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats = -4122
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Boolean "outlier" flags for each algorithm's imputed value, rows 2-18.
$values = @{
    2  = @(0,0,0)
    3  = @(0,0,0)
    4  = @(1,0,0)
    5  = @(0,0,0)
    6  = @(0,0,0)
    7  = @(0,0,0)
    8  = @(0,0,0)
    9  = @(1,1,0)
    10 = @(0,0,0)
    11 = @(0,0,0)
    12 = @(0,0,0)
    13 = @(0,0,0)
    14 = @(1,0,0)
    15 = @(0,0,0)
    16 = @(0,0,0)
    17 = @(0,0,0)
    18 = @(0,0,0)
}

foreach ($row in $values.Keys) {
    $triple = $values[$row]
    $ws.Cells.Item($row, 6).Value = [bool]$triple[0]
    $ws.Cells.Item($row, 7).Value = [bool]$triple[1]
    $ws.Cells.Item($row, 8).Value = [bool]$triple[2]
}
